# [fix] Battle Protocol design
#
# The "Battle" resource's create/delete operations (POST/DELETE on
# /Sprout/Battles) and the "win/lose" update operation (PUT on
# /Sprout/Battles/{battleId}/Teams used for results) are removed from the
# protocol design. Their usage notes, parameters and sample
# responses are cleared and the now-unavailable HTTP methods are marked
# with a strikethrough, matching the styling already used elsewhere in the
# sheet for methods that a resource does not support.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Battles resource (rows 2-5): POST (create battle) and DELETE
#     (delete battle) are no longer supported -> clear their usage /
#     parameter / response cells and strike the method name through.

# Row 3 : POST /Sprout/Battles ("バトルの作成")
$ws.Range("B3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("C3").Font.Strikethrough = $true
$ws.Rows.Item(3).AutoFit()

# Row 5 : DELETE /Sprout/Battles ("バトルの削除")
$ws.Range("B5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("C5").Font.Strikethrough = $true

# --- Teams resource (rows 10-13): DELETE is (now) unsupported as well.
$ws.Range("C13").Font.Strikethrough = $true

# --- Battle results (rows 14-17): PUT ("勝敗情報の更新") is removed ->
#     clear its usage / parameter cells and strike PUT through.
$ws.Range("B16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("C16").Font.Strikethrough = $true

# --- View state: selection moved to F14, scrolled so row 10 is at the top.
$null = $ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
